# Apply changes: rename sheet, add new rows, set column widths, update selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "User Data" to "Item Data"
$ws.Name = "Item Data"

# Add new rows for items (Item ID, Item Name="test", Price)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "test"
$ws.Cells.Item(5, 3).Value = 1

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "test"
$ws.Cells.Item(6, 3).Value = 1

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "test"
$ws.Cells.Item(7, 3).Value = 1

# Autofit columns (bestFit / customWidth)
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()

# Update selection to row 5 (entire row)
$ws.Rows.Item(5).Select() | Out-Null
